# Corrects IFRS company_list figures for rows 2-9 (sheet "company_list"):
# - Rows 2-6: replace mis-scaled financial figures with the corrected values;
#   a couple of ratio cells (AD/AH) are removed entirely where no longer reported.
# - Rows 7-9: clear all financial data columns, keeping only id/name columns (A-C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$values = @{
    "D" = 2545;
    "E" = 82;
    "F" = 82;
    "G" = 81;
    "H" = 55;
    "I" = 56;
    "J" = -1;
    "K" = 1929;
    "L" = 1488;
    "M" = 441;
    "N" = 429;
    "O" = 12;
    "P" = 49;
    "Q" = 129;
    "R" = -163;
    "S" = 32;
    "T" = 181;
    "U" = -51;
    "V" = 858;
    "W" = 3.21;
    "X" = 2.18;
    "Y" = 15.02;
    "Z" = 3.07;
    "AA" = 336.98;
    "AB" = 453.11;
    "AC" = 567;
    "AE" = 4348;
    "AF" = 0;
    "AG" = 0;
    "AI" = 0;
    "AJ" = 9868460
}
foreach ($col in $values.Keys) {
    $ws.Range("${col}2").Value = $values[$col]
}
@("AD2","AH2") | ForEach-Object { $ws.Range($_).ClearContents() }

# --- Row 3 ---
$values = @{
    "D" = 2778;
    "E" = 134;
    "F" = 134;
    "G" = 93;
    "H" = 66;
    "I" = 66;
    "J" = 0;
    "K" = 2475;
    "L" = 1968;
    "M" = 507;
    "N" = 505;
    "O" = 2;
    "P" = 49;
    "Q" = 116;
    "R" = -371;
    "S" = 295;
    "T" = 408;
    "U" = -292;
    "V" = 1179;
    "W" = 4.81;
    "X" = 2.39;
    "Y" = 14.23;
    "Z" = 3.01;
    "AA" = 388.43;
    "AB" = 891.9;
    "AC" = 673;
    "AE" = 5113;
    "AF" = 0;
    "AG" = 0;
    "AI" = 0;
    "AJ" = 9868460
}
foreach ($col in $values.Keys) {
    $ws.Range("${col}3").Value = $values[$col]
}
@("AD3","AH3") | ForEach-Object { $ws.Range($_).ClearContents() }

# --- Row 4 ---
$values = @{
    "D" = 2980;
    "E" = 155;
    "F" = 155;
    "G" = 119;
    "H" = 92;
    "I" = 92;
    "J" = 1;
    "K" = 2851;
    "L" = 2237;
    "M" = 615;
    "N" = 601;
    "O" = 13;
    "P" = 49;
    "Q" = 44;
    "R" = -256;
    "S" = 214;
    "T" = 308;
    "U" = -265;
    "V" = 1463;
    "W" = 5.21;
    "X" = 3.09;
    "Y" = 16.56;
    "Z" = 3.46;
    "AA" = 363.9;
    "AB" = 1079.92;
    "AC" = 928;
    "AE" = 6094;
    "AF" = 0;
    "AG" = 0;
    "AI" = 0;
    "AJ" = 9868460
}
foreach ($col in $values.Keys) {
    $ws.Range("${col}4").Value = $values[$col]
}
@("AD4","AH4") | ForEach-Object { $ws.Range($_).ClearContents() }

# --- Row 5 ---
$values = @{
    "D" = 2942;
    "E" = 190;
    "F" = 190;
    "G" = 137;
    "H" = 94;
    "I" = 92;
    "J" = 2;
    "K" = 2986;
    "L" = 2123;
    "M" = 863;
    "N" = 847;
    "O" = 15;
    "P" = 66;
    "Q" = 181;
    "R" = -265;
    "S" = 105;
    "T" = 160;
    "U" = 21;
    "V" = 1369;
    "W" = 6.47;
    "X" = 3.2;
    "Y" = 12.68;
    "Z" = 3.22;
    "AA" = 245.97;
    "AB" = 1248.19;
    "AC" = 905;
    "AD" = 4.53;
    "AE" = 6436;
    "AF" = 0.64;
    "AG" = 100;
    "AH" = 2.44;
    "AI" = 14.33;
    "AJ" = 13168460
}
foreach ($col in $values.Keys) {
    $ws.Range("${col}5").Value = $values[$col]
}

# --- Row 6 ---
$values = @{
    "D" = 3322;
    "E" = 174;
    "F" = 174;
    "G" = 144;
    "H" = 117;
    "I" = 115;
    "K" = 3176;
    "L" = 2221;
    "M" = 955;
    "N" = 937;
    "P" = 66;
    "Q" = 186;
    "R" = -234;
    "S" = 44;
    "T" = 252;
    "U" = -66;
    "V" = 1446;
    "W" = 5.25;
    "X" = 3.54;
    "Y" = 12.88;
    "Z" = 3.81;
    "AA" = 232.53;
    "AB" = 1377.78;
    "AC" = 873;
    "AD" = 5.6;
    "AE" = 7246;
    "AF" = 0.67;
    "AG" = 120;
    "AH" = 2.46;
    "AI" = 13.51;
    "AJ" = 13168460
}
foreach ($col in $values.Keys) {
    $ws.Range("${col}6").Value = $values[$col]
}

# --- Row 7 ---
@("D7","E7","G7","H7","I7","K7","L7","M7","N7","P7","Q7","R7","S7","T7","U7","W7","X7","Y7","Z7","AA7","AC7","AD7","AE7","AF7","AG7","AH7","AI7") | ForEach-Object { $ws.Range($_).ClearContents() }

# --- Row 8 ---
@("D8","E8","G8","H8","I8","K8","L8","M8","N8","P8","Q8","R8","S8","T8","U8","W8","X8","Y8","Z8","AA8","AC8","AD8","AE8","AF8","AG8","AH8","AI8") | ForEach-Object { $ws.Range($_).ClearContents() }

# --- Row 9 ---
@("D9","E9","G9","H9","I9","K9","L9","M9","N9","P9","Q9","R9","S9","T9","U9","W9","X9","Y9","Z9","AA9","AC9","AD9","AE9","AF9","AG9","AH9","AI9") | ForEach-Object { $ws.Range($_).ClearContents() }
